# Update cryptos price (D) and volume-change (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.263.10"
$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("D3").Value = "'1.823.14"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  -1.15%  "
$ws.Range("D5").Value = "'314.53"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  -1.21%  "
$ws.Range("D7").Value = "'0.4279"
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("D8").Value = "'0.3684"
$ws.Range("E8").Value = "  -2.60%  "
$ws.Range("D9").Value = "'0.07241"
$ws.Range("E9").Value = "  -2.58%  "
$ws.Range("D10").Value = "'0.8626"
$ws.Range("E10").Value = "  -2.38%  "
$ws.Range("D11").Value = "'21.01"
$ws.Range("E11").Value = "  -2.49%  "
$ws.Range("D12").Value = "'1.817.80"
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("D13").Value = "'6.676"
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("D14").Value = "'0.07109"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "'5.313"
$ws.Range("E15").Value = "  -3.17%  "
$ws.Range("D16").Value = "'89.11"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").Value = "'1.007"
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").Value = "'0.000008879"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("D20").Value = "'15.06"
$ws.Range("E20").Value = "  -2.59%  "
$ws.Range("D21").Value = "'27.281.19"
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("D22").Value = "'5.149"
$ws.Range("E22").Value = "  -2.32%  "
$ws.Range("E23").Value = "  -2.32%  "
$ws.Range("D24").Value = "'2.048.68"
$ws.Range("E24").Value = "  -1.74%  "
$ws.Range("D25").Value = "'2.007"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("D26").Value = "'153.49"
$ws.Range("E26").Value = "  -2.12%  "
$ws.Range("D27").Value = "'18.40"
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("E28").Value = "  +6.63%  "
$ws.Range("D29").Value = "'5.240"
$ws.Range("E29").Value = "  -3.64%  "
$ws.Range("D30").Value = "'116.37"
$ws.Range("E30").Value = "  -3.53%  "
$ws.Range("D31").Value = "'0.08905"
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("D32").Value = "'1.203"
$ws.Range("E32").Value = "  -1.81%  "
$ws.Range("D33").Value = "'0.7605"
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("D34").Value = "'4.464"
$ws.Range("E34").Value = "  -2.10%  "
$ws.Range("D35").Value = "'2.812"
$ws.Range("E35").Value = "  -7.30%  "
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("D37").Value = "'1.118"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("D38").Value = "'0.01972"
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("D39").Value = "'0.05288"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").Value = "'2.911"
$ws.Range("E40").Value = "  +1.77%  "
$ws.Range("D41").Value = "'7.139"
$ws.Range("E41").Value = "  +2.64%  "
$ws.Range("D42").Value = "'0.1686"
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").Value = "'0.5046"
$ws.Range("E43").Value = "  -2.69%  "
$ws.Range("D44").Value = "'8.642"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("D45").Value = "'10.60"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("D46").Value = "'106.74"
$ws.Range("E46").Value = "  -2.99%  "
$ws.Range("D47").Value = "'0.4758"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("D48").Value = "'1.005"
$ws.Range("E48").Value = "  -1.39%  "
$ws.Range("D49").Value = "'0.06401"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").Value = "'1.664"
$ws.Range("E50").Value = "  -2.83%  "
$ws.Range("D51").Value = "'1.814"
$ws.Range("E51").Value = "  -2.06%  "
